# Update "manaCost" (column D) values on the WeaponCombatList sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WeaponCombatList")

$ws.Range("D2").Value  = 10
$ws.Range("D3").Value  = 10
$ws.Range("D4").Value  = 15
$ws.Range("D5").Value  = 10
$ws.Range("D6").Value  = 10
$ws.Range("D7").Value  = 15
$ws.Range("D8").Value  = 10
$ws.Range("D9").Value  = 10
$ws.Range("D10").Value = 15
$ws.Range("D11").Value = 10
$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 15
$ws.Range("D14").Value = 40
$ws.Range("D15").Value = 25
$ws.Range("D16").Value = 55

# Make WeaponCombatList the active sheet/tab (was StatBoostEffectList)
# and update its selected cell.
$ws.Activate()
$ws.Range("D17").Select()
